$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2369
$ws1.Range("F7").Value = 8205
$ws1.Range("F11").Value = 4705
$ws1.Range("F14").Value = 837
$ws1.Range("F16").Value = 1324
$ws1.Range("F19").Value = 6636
$ws1.Range("F23").Value = 4460
$ws1.Range("F24").Value = 340
$ws1.Range("F25").Value = 738
$ws1.Range("F26").Value = 2125
$ws1.Range("F28").Value = 376
$ws1.Range("F30").Value = 116
$ws1.Range("F33").Value = 111
$ws1.Range("F34").Value = 348
$ws1.Range("F40").Value = 1266
$ws1.Range("F44").Value = 1264

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 411
$ws2.Range("F12").Value = 215

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F11").Value = 1600
$ws3.Range("F12").Value = 1928
$ws3.Range("F13").Value = 409
$ws3.Range("F14").Value = 308

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2369
$ws4.Range("F11").Value = 1600
$ws4.Range("F13").Value = 1928
$ws4.Range("F14").Value = 4705
$ws4.Range("F15").Value = 411
$ws4.Range("F17").Value = 837
$ws4.Range("F20").Value = 1324
$ws4.Range("F23").Value = 6636
$ws4.Range("F25").Value = 308
$ws4.Range("F27").Value = 340
$ws4.Range("F28").Value = 2125
$ws4.Range("F30").Value = 376
$ws4.Range("F31").Value = 116
$ws4.Range("F33").Value = 215
$ws4.Range("F35").Value = 111
$ws4.Range("F36").Value = 348
$ws4.Range("F40").Value = 1266
$ws4.Range("F46").Value = 1264
